# Update specific imputed values in Sheet1 (RandomForest result data)
# as produced by a re-run of the imputation algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.0528
$ws.Range("C4").Value = -12.3008
$ws.Range("A11").Value = -21.5203
$ws.Range("A12").Value = -21.561
$ws.Range("C14").Value = -13.13939999999999
$ws.Range("A15").Value = -21.9418
$ws.Range("C26").Value = -12.34780000000001
$ws.Range("A27").Value = -21.642
$ws.Range("A28").Value = -21.8323
$ws.Range("A31").Value = -21.56160000000001
$ws.Range("C31").Value = -12.4642
$ws.Range("A32").Value = -21.32619999999999
$ws.Range("C35").Value = -13.55250000000001
$ws.Range("A36").Value = -21.05420000000001
$ws.Range("C37").Value = -14.1649
$ws.Range("A38").Value = -19.57389999999999
$ws.Range("C39").Value = -13.39170000000001
$ws.Range("C40").Value = -13.3474
$ws.Range("C45").Value = -14.2281
$ws.Range("A46").Value = -21.44699999999998
$ws.Range("C52").Value = -11.1538
$ws.Range("A54").Value = -21.6162
$ws.Range("A55").Value = -22.21330000000001
$ws.Range("A56").Value = -22.15030000000002
$ws.Range("C57").Value = -14.5075
$ws.Range("A67").Value = -21.57249999999998
$ws.Range("A69").Value = -21.69529999999997
$ws.Range("A72").Value = -21.59409999999999
$ws.Range("A73").Value = -19.97669999999999
$ws.Range("C81").Value = -12.58060000000001
$ws.Range("A83").Value = -21.848
$ws.Range("C83").Value = -12.61019999999999
$ws.Range("A86").Value = -22.30980000000001
$ws.Range("A91").Value = -21.42180000000001
$ws.Range("A93").Value = -21.4482
$ws.Range("A99").Value = -20.35939999999999
$ws.Range("C100").Value = -12.8923
$ws.Range("C102").Value = -13.3719
